# Update column C (rows 2-92) from serial date 45188 (2023-09-19) to 45189 (2023-09-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($r = 2; $r -le 92; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
